$d = $word.ActiveDocument

# --- Title paragraph: add tab stop, first-line indent, and move the
#     "_GoBack" bookmark here (Word keeps a single _GoBack bookmark at the
#     most-recent edit location; adding a new one removes the old one). ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.ParagraphFormat.TabStops.Add(252)
$p1.Range.ParagraphFormat.FirstLineIndent = 14.4

# Position 0 (the very start of the document) has a quirky interaction when
# used directly as a zero-length bookmark range, so nudge everything over by
# inserting a placeholder character, bookmark at the now-safe offset, then
# remove the placeholder again.
$r0 = $d.Range(0, 0)
$r0.InsertBefore("X")
$rb = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $rb)
$rdel = $d.Range(0, 1)
$rdel.Delete()

# --- Abstract paragraph: add a tab stop. ---
$pAbstract = $d.Paragraphs.Item(5)
$pAbstract.Range.ParagraphFormat.TabStops.Add(324)

Write-Host "Done"
